# Update "paises" (countries) & "provincias Spain" data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---------------------------------------------------------------------------
# Update the "last updated" timestamp shown in A1
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 10:23"

# ---------------------------------------------------------------------------
# Helper: update a full data row (columns B..H) for a given sheet row number
# ---------------------------------------------------------------------------
function Set-Row($row, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 2).Value = $casosTotales
    $ws.Cells.Item($row, 3).Value = $nuevosCasos
    $ws.Cells.Item($row, 4).Value = $casosActivos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $casosCriticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Row 7  - Rusia
Set-Row 7 897599 4945 703175 179293 0 130 15131

# Row 25 - Filipinas
Set-Row 25 139538 2987 68432 68794 0 19 2312

# Row 47 - Singapur
Set-Row 47 55353 61 49609 5717 0 0 27

# Row 49 - Polonia
Set-Row 49 52410 0 37150 13451 0 0 1809

# Row 108 - Hungria
Set-Row 108 4746 15 3527 614 0 0 605

# Row 111 - Hong Kong
Set-Row 111 4149 0 2917 1174 0 3 58

# Row 122 - Sri Lanka
Set-Row 122 2871 0 2622 238 0 0 11

# Row 128 - Lituania
Set-Row 128 2283 18 1679 523 0 0 81

# Row 130 - Estonia
Set-Row 130 2167 9 1968 136 0 0 63

# ---------------------------------------------------------------------------
# Montserrat moved above Islas Malvinas in the source country list, so the
# data that used to sit in row 213 (Islas Malvinas) now belongs to row 214,
# and vice versa - swap both the country label and the statistics.
# ---------------------------------------------------------------------------
$ws.Cells.Item(213, 1).Value = "Montserrat"
Set-Row 213 13 0 12 0 0 0 1

$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
Set-Row 214 13 0 13 0 0 0 0
